# Updated cryptos list with GitHub Actions - refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ D = "26.205.55";  E = "  -3.91%  " }
    3  = @{ D = "1.659.74";   E = "  -2.49%  " }
    4  = @{ E = "  +0.23%  " }
    5  = @{ E = "  -2.31%  " }
    6  = @{ D = "0.5151";     E = "  -2.85%  " }
    7  = @{ E = "  +0.26%  " }
    8  = @{ D = "0.2581";     E = "  -2.89%  " }
    9  = @{ D = "0.06453";    E = "  -1.93%  " }
    10 = @{ D = "20.00";      E = "  -3.36%  " }
    11 = @{ D = "0.07807" }
    12 = @{ D = "1.656.79";   E = "  -2.71%  " }
    13 = @{ E = "  -4.18%  " }
    14 = @{ D = "1.887.29";   E = "  -2.52%  " }
    15 = @{ D = "0.5551";     E = "  -3.95%  " }
    16 = @{ D = "0.0$([char]0x2085)8079"; E = "  -0.78%  " }
    17 = @{ D = "64.36";      E = "  -4.48%  " }
    18 = @{ D = "26.222.92";  E = "  -3.81%  " }
    19 = @{ D = "211.68";     E = "  -1.68%  " }
    20 = @{ E = "  +0.24%  " }
    21 = @{ D = "4.439";      E = "  -3.70%  " }
    22 = @{ D = "10.07";      E = "  -2.51%  " }
    23 = @{ D = "5.992";      E = "  +0.45%  " }
    24 = @{ E = "  +0.23%  " }
    25 = @{ D = "144.32";     E = "  +0.24%  " }
    26 = @{ E = "  +2.98%  " }
    27 = @{ D = "0.1169";     E = "  -2.30%  " }
    28 = @{ D = "6.989";      E = "  -2.91%  " }
    29 = @{ D = "15.82";      E = "  -1.85%  " }
    30 = @{ D = "0.05214";    E = "  -2.75%  " }
    31 = @{ D = "1.255";      E = "  -2.25%  " }
    32 = @{ D = "3.366";      E = "  -2.76%  " }
    33 = @{ D = "3.226";      E = "  -5.07%  " }
    34 = @{ D = "1.574";      E = "  -3.91%  " }
    35 = @{ D = "2.762";      E = "  -3.52%  " }
    36 = @{ D = "0.9331";     E = "  -1.33%  " }
    37 = @{ D = "2.371";      E = "  -1.90%  " }
    38 = @{ D = "1.176.84";   E = "  +13.21%  " }
    39 = @{ D = "0.5704";     E = "  -1.87%  " }
    40 = @{ D = "0.01598";    E = "  -1.84%  " }
    41 = @{ D = "0.8464";     E = "  +0.90%  " }
    42 = @{ E = "  +0.20%  " }
    43 = @{ D = "5.685";      E = "  -1.53%  " }
    44 = @{ D = "100.69";     E = "  -0.24%  " }
    45 = @{ D = "1.797.28";   E = "  -2.54%  " }
    46 = @{ E = "  -1.67%  " }
    47 = @{ E = "  +0.45%  " }
    48 = @{ D = "56.04";      E = "  -3.00%  " }
    49 = @{ D = "1.005";      E = "  +0.04%  " }
    50 = @{ D = "7.872";      E = "  -2.31%  " }
    51 = @{ D = "0.05065";    E = "  -3.10%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        $cell = $ws.Cells.Item($row, 4)
        # Column D holds text-formatted price strings (e.g. "1.659.74"); force
        # text so Excel doesn't reinterpret/renormalize them as numbers.
        $cell.NumberFormat = "@"
        $cell.Value = $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        # Column E values are space-padded percentages (e.g. "  -3.91%  "),
        # which Excel never reinterprets as numbers, so no format coercion
        # is needed here.
        $ws.Cells.Item($row, 5).Value = $vals["E"]
    }
}
